$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 563.9487
$ws.Range("J17").Value = 573.5263
$ws.Range("L17").Value = 1720.5789
$ws.Range("N17").Value = -2056.5789

$ws.Range("H43").Value = 938.2143
$ws.Range("I43").Value = 720.5
$ws.Range("J43").Value = 974.5
$ws.Range("K43").Value = 720.5
$ws.Range("L43").Value = 974.5
$ws.Range("N43").Value = -1112.5
$ws.Range("M43").Value = -651.5

$ws.Range("H70").Value = 1646.4
$ws.Range("I70").Value = 1126
$ws.Range("J70").Value = 1993.3334
$ws.Range("K70").Value = 3378
$ws.Range("L70").Value = 5980.0002
$ws.Range("M70").Value = -3108
$ws.Range("N70").Value = -6520.0002

$ws.Range("H73").Value = 1646.4
$ws.Range("I73").Value = 1126
$ws.Range("J73").Value = 1993.3334
$ws.Range("K73").Value = 3378
$ws.Range("L73").Value = 5980.0002
$ws.Range("M73").Value = -2442
$ws.Range("N73").Value = -7852.0002

$ws.Range("H129").Value = 1175.4524
$ws.Range("J129").Value = 1389.875
$ws.Range("L129").Value = 4169.625
$ws.Range("N129").Value = -14169.625

$ws.Range("H137").Value = 7896900
$ws.Range("I137").Value = 3847792.5
$ws.Range("J137").Value = 16669967
$ws.Range("K137").Value = 11543377.5
$ws.Range("L137").Value = 50009901
$ws.Range("M137").Value = -11540827.5
$ws.Range("N137").Value = -50015001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1779.6786
$ws.Range("I61").Value = 1397.3462
$ws.Range("J61").Value = 6750
$ws.Range("K61").Value = 1397.3462
$ws.Range("L61").Value = 6750
$ws.Range("M61").Value = -1185.3462
$ws.Range("N61").Value = -7174

$ws.Range("H136").Value = 1779.6786
$ws.Range("I136").Value = 1397.3462
$ws.Range("J136").Value = 6750
$ws.Range("K136").Value = 4192.0386
$ws.Range("L136").Value = 20250
$ws.Range("M136").Value = -1642.0386
$ws.Range("N136").Value = -25350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4984.5
$ws.Range("I107").Value = 5336.75
$ws.Range("J107").Value = 3751.625
$ws.Range("K107").Value = 5336.75
$ws.Range("L107").Value = 3751.625
$ws.Range("M107").Value = -3416.75
$ws.Range("N107").Value = -7591.625

$ws.Range("H113").Value = 3340
$ws.Range("I113").Value = 3340
$ws.Range("K113").Value = 3340
$ws.Range("M113").Value = -1170

$ws.Range("H134").Value = 68461.53
$ws.Range("I134").Value = 87439.17
$ws.Range("J134").Value = 2039.8
$ws.Range("K134").Value = 262317.51
$ws.Range("L134").Value = 6119.4
$ws.Range("M134").Value = -259782.51
$ws.Range("N134").Value = -11189.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 18141.646
$ws.Range("J4").Value = 31167.334
$ws.Range("L4").Value = 31167.334
$ws.Range("N4").Value = -31391.334

$ws.Range("H31").Value = 2299.8333
$ws.Range("I31").Value = 1694.75
$ws.Range("K31").Value = 1694.75
$ws.Range("M31").Value = -1399.75

$ws.Range("H34").Value = 2299.8333
$ws.Range("I34").Value = 1694.75
$ws.Range("K34").Value = 1694.75
$ws.Range("M34").Value = -1492.75

$ws.Range("H58").Value = 1227.0513
$ws.Range("I58").Value = 1100.2963
$ws.Range("J58").Value = 1512.25
$ws.Range("K58").Value = 1100.2963
$ws.Range("L58").Value = 1512.25
$ws.Range("M58").Value = -897.2963
$ws.Range("N58").Value = -1918.25

$ws.Range("H132").Value = 3755.6924
$ws.Range("I132").Value = 3338.1
$ws.Range("J132").Value = 5147.6665
$ws.Range("K132").Value = 10014.3
$ws.Range("L132").Value = 15442.9995
$ws.Range("M132").Value = -7484.299999999999
$ws.Range("N132").Value = -20502.9995

$ws.Range("H134").Value = 5903.222
$ws.Range("I134").Value = 8533.4375
$ws.Range("J134").Value = 2077.4546
$ws.Range("K134").Value = 25600.3125
$ws.Range("L134").Value = 6232.3638
$ws.Range("M134").Value = -23065.3125
$ws.Range("N134").Value = -11302.3638

$ws.Range("H136").Value = 1227.0513
$ws.Range("I136").Value = 1100.2963
$ws.Range("J136").Value = 1512.25
$ws.Range("K136").Value = 3300.8889
$ws.Range("L136").Value = 4536.75
$ws.Range("M136").Value = -750.8888999999999
$ws.Range("N136").Value = -9636.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 823354.25
$ws.Range("I132").Value = 1197018
$ws.Range("J132").Value = 1294
$ws.Range("K132").Value = 10773162
$ws.Range("L132").Value = 11646
$ws.Range("M132").Value = -10770632
$ws.Range("N132").Value = -16706

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1936.1904
$ws.Range("I102").Value = 1758
$ws.Range("J102").Value = 2381.6667
$ws.Range("K102").Value = 1758
$ws.Range("L102").Value = 2381.6667
$ws.Range("M102").Value = -136
$ws.Range("N102").Value = -5625.6667

$ws.Range("H132").Value = 2844
$ws.Range("I132").Value = 2297.4
$ws.Range("K132").Value = 6892.200000000001
$ws.Range("M132").Value = -4362.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2950
$ws.Range("J2").Value = 5600
$ws.Range("L2").Value = 5600
$ws.Range("N2").Value = -5824

$ws.Range("H132").Value = 2444.6365
$ws.Range("I132").Value = 1971.2727
$ws.Range("K132").Value = 5913.8181
$ws.Range("M132").Value = -3383.8181

$ws.Range("H136").Value = 1525.3103
$ws.Range("I136").Value = 1267
$ws.Range("K136").Value = 3801
$ws.Range("M136").Value = -1251

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 27513012
$ws.Range("I5").Value = 50
$ws.Range("J5").Value = 36684000
$ws.Range("K5").Value = 50
$ws.Range("L5").Value = 36684000
$ws.Range("M5").Value = 62
$ws.Range("N5").Value = -36684224

$ws.Range("H81").Value = 5484.75
$ws.Range("I81").Value = 9693.817999999999
$ws.Range("J81").Value = 1923.2307
$ws.Range("K81").Value = 19387.636
$ws.Range("L81").Value = 3846.4614
$ws.Range("M81").Value = -18326.636
$ws.Range("N81").Value = -5968.4614

$ws.Range("H84").Value = 5484.75
$ws.Range("I84").Value = 9693.817999999999
$ws.Range("J84").Value = 1923.2307
$ws.Range("K84").Value = 96938.17999999999
$ws.Range("L84").Value = 19232.307
$ws.Range("M84").Value = -91634.17999999999
$ws.Range("N84").Value = -29840.307

$ws.Range("H107").Value = 1440.9131
$ws.Range("I107").Value = 1095.5714
$ws.Range("J107").Value = 1978.1111
$ws.Range("K107").Value = 3286.7142
$ws.Range("L107").Value = 5934.3333
$ws.Range("M107").Value = -1366.7142
$ws.Range("N107").Value = -9774.3333

$ws.Range("H118").Value = 43429.332
$ws.Range("J118").Value = 43429.332
$ws.Range("L118").Value = 43429.332
$ws.Range("N118").Value = -46743.332

$ws.Range("H132").Value = 3012.4856
$ws.Range("I132").Value = 2598.0386
$ws.Range("J132").Value = 4209.778
$ws.Range("K132").Value = 7794.1158
$ws.Range("L132").Value = 12629.334
$ws.Range("M132").Value = -5264.1158
$ws.Range("N132").Value = -17689.334

$ws.Range("H136").Value = 1473.2368
$ws.Range("I136").Value = 1222.3055
$ws.Range("J136").Value = 5990
$ws.Range("K136").Value = 3666.9165
$ws.Range("L136").Value = 17970
$ws.Range("M136").Value = -1116.9165
$ws.Range("N136").Value = -23070
